$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# $apos is a literal apostrophe. Prefixing a value with it forces Excel to store
# the cell as literal text instead of auto-converting a numeric-looking string to
# a floating point number (which would corrupt formatting such as trailing zeros,
# e.g. "19.50" -> 19.5, or "0.0000169" -> 1.69E-05).
$apos = [string][char]39

# ---- Price (D) and Volume(1h) (E) updates ----
$ws.Range('D2').Value = '60.839.18'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.400.10'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.58%  '
$ws.Range('D5').Value = $apos + '565.87'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = $apos + '141.81'
$ws.Range('E6').Value = '  +2.64%  '
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('E8').Value = '  +2.56%  '
$ws.Range('D9').Value = '2.407.26'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('E12').Value = '  +2.88%  '
$ws.Range('E13').Value = '  +2.89%  '
$ws.Range('D14').Value = $apos + '26.42'
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').Value = $apos + '0.0000169'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '2.837.00'
$ws.Range('D17').Value = '60.623.02'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '2.408.42'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = $apos + '8.05'
$ws.Range('E19').Value = '  +3.29%  '
$ws.Range('D20').Value = $apos + '10.69'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').Value = $apos + '324.17'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('D23').Value = $apos + '6.06'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('E25').Value = '  +5.01%  '
$ws.Range('D26').Value = $apos + '65.11'
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('D27').Value = $apos + '583.55'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').Value = $apos + '8.22'
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('D29').Value = '0.0₃0940'
$ws.Range('E29').Value = '  +2.91%  '
$ws.Range('D30').Value = '2.518.46'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('E31').Value = '  +3.05%  '
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('E35').Value = '  +5.84%  '
$ws.Range('D37').Value = $apos + '153.45'
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('E38').Value = '  +1.88%  '
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('D40').Value = $apos + '18.30'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').Value = $apos + '5.18'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').Value = $apos + '2.51'
$ws.Range('E43').Value = '  +12.16%  '
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('D45').Value = $apos + '41.61'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').Value = '0.0₆0279'
$ws.Range('E46').Value = '  +6.69%  '
$ws.Range('D47').Value = $apos + '141.47'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').Value = $apos + '19.50'
$ws.Range('E51').Value = '  +1.79%  '

# Reset style back to "Normal" on the text-forced cells above so that no extra
# number-format/style gets attached to them (keeps styling identical to original).
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D51').Style = "Normal"
